# Auto-generated Excel COM-interop script to apply Tonberry_Profits.xlsx value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 779.1111
$ws.Range("I12").Value = 702.3333
$ws.Range("K12").Value = 702.3333
$ws.Range("M12").Value = -532.3333

$ws.Range("H62").Value = 2700.6667
$ws.Range("I62").Value = 2431.6667
$ws.Range("K62").Value = 2431.6667
$ws.Range("M62").Value = -1807.6667

$ws.Range("H65").Value = 2700.6667
$ws.Range("I65").Value = 2431.6667
$ws.Range("K65").Value = 12158.3335
$ws.Range("M65").Value = -9038.333500000001

$ws.Range("H70").Value = 21090.908
$ws.Range("I70").Value = 13783.333
$ws.Range("J70").Value = 29860
$ws.Range("K70").Value = 41349.999
$ws.Range("L70").Value = 89580
$ws.Range("M70").Value = -41079.999
$ws.Range("N70").Value = -90120

$ws.Range("H73").Value = 21090.908
$ws.Range("I73").Value = 13783.333
$ws.Range("J73").Value = 29860
$ws.Range("K73").Value = 41349.999
$ws.Range("L73").Value = 89580
$ws.Range("M73").Value = -40413.999
$ws.Range("N73").Value = -91452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2423.6667
$ws.Range("I32").Value = 1669.7858
$ws.Range("J32").Value = 6645.4
$ws.Range("K32").Value = 1669.7858
$ws.Range("L32").Value = 6645.4
$ws.Range("M32").Value = -1382.7858
$ws.Range("N32").Value = -7219.4

$ws.Range("H45").Value = 1400.5186
$ws.Range("I45").Value = 890.1579
$ws.Range("J45").Value = 2612.625
$ws.Range("K45").Value = 890.1579
$ws.Range("L45").Value = 2612.625
$ws.Range("M45").Value = -513.1579
$ws.Range("N45").Value = -3366.625

$ws.Range("H61").Value = 4372
$ws.Range("I61").Value = 3744
$ws.Range("K61").Value = 3744
$ws.Range("M61").Value = -3532

$ws.Range("H97").Value = 525.875
$ws.Range("I97").Value = 525.875
$ws.Range("K97").Value = 525.875
$ws.Range("M97").Value = -29.875

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 1704
$ws.Range("I132").Value = 1350.7693
$ws.Range("K132").Value = 4052.3079
$ws.Range("M132").Value = -1522.3079

$ws.Range("H136").Value = 4372
$ws.Range("I136").Value = 3744
$ws.Range("K136").Value = 11232
$ws.Range("M136").Value = -8682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27499.75
$ws.Range("I82").Value = 21666.334
$ws.Range("K82").Value = 21666.334
$ws.Range("M82").Value = -21283.334

$ws.Range("H85").Value = 27499.75
$ws.Range("I85").Value = 21666.334
$ws.Range("K85").Value = 21666.334
$ws.Range("M85").Value = -20340.334

$ws.Range("H134").Value = 12694.5
$ws.Range("I134").Value = 14438.583
$ws.Range("J134").Value = 7462.25
$ws.Range("K134").Value = 43315.749
$ws.Range("L134").Value = 22386.75
$ws.Range("M134").Value = -40780.749
$ws.Range("N134").Value = -27456.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 34319.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 34319.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 34319.5
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -34799.5

$ws.Range("H27").Value = 34319.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 34319.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 34319.5
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -34703.5

$ws.Range("H122").Value = 1601.7333
$ws.Range("I122").Value = 1412.6316
$ws.Range("K122").Value = 4237.8948
$ws.Range("M122").Value = -1787.8948

$ws.Range("H132").Value = 1835.5834
$ws.Range("I132").Value = 1171.95
$ws.Range("J132").Value = 5153.75
$ws.Range("K132").Value = 3515.85
$ws.Range("L132").Value = 15461.25
$ws.Range("M132").Value = -985.8500000000004
$ws.Range("N132").Value = -20521.25

$ws.Range("H134").Value = 916.44446
$ws.Range("I134").Value = 873.9375
$ws.Range("K134").Value = 2621.8125
$ws.Range("M134").Value = -86.8125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 606.8570999999999
$ws.Range("I5").Value = 517.4545000000001
$ws.Range("K5").Value = 1552.3635
$ws.Range("M5").Value = -1440.3635

$ws.Range("H26").Value = 261.5
$ws.Range("J26").Value = 248.66667
$ws.Range("L26").Value = 746.00001
$ws.Range("N26").Value = -1322.00001

$ws.Range("H98").Value = 732.3333
$ws.Range("J98").Value = 797
$ws.Range("L98").Value = 2391
$ws.Range("N98").Value = -5387

$ws.Range("H114").Value = 1820.75
$ws.Range("I114").Value = 386.66666
$ws.Range("J114").Value = 2681.2
$ws.Range("K114").Value = 1159.99998
$ws.Range("L114").Value = 8043.599999999999
$ws.Range("M114").Value = 2094.00002
$ws.Range("N114").Value = -14551.6

$ws.Range("H132").Value = 1477.2941
$ws.Range("I132").Value = 1481.5454
$ws.Range("K132").Value = 13333.9086
$ws.Range("M132").Value = -10803.9086

$ws.Range("H135").Value = 606.8570999999999
$ws.Range("I135").Value = 517.4545000000001
$ws.Range("K135").Value = 4657.0905
$ws.Range("M135").Value = -2122.0905

$ws.Range("H140").Value = 1753.9117
$ws.Range("I140").Value = 778.75
$ws.Range("J140").Value = 2053.9614
$ws.Range("K140").Value = 2336.25
$ws.Range("L140").Value = 6161.8842
$ws.Range("M140").Value = 2843.75
$ws.Range("N140").Value = -16521.8842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5001.75
$ws.Range("I70").Value = 5302.8
$ws.Range("K70").Value = 5302.8
$ws.Range("M70").Value = -5032.8

$ws.Range("H73").Value = 5001.75
$ws.Range("I73").Value = 5302.8
$ws.Range("K73").Value = 5302.8
$ws.Range("M73").Value = -4366.8

$ws.Range("H97").Value = 1478.7273
$ws.Range("I97").Value = 453.8
$ws.Range("J97").Value = 2332.8333
$ws.Range("K97").Value = 453.8
$ws.Range("L97").Value = 2332.8333
$ws.Range("M97").Value = 42.19999999999999
$ws.Range("N97").Value = -3324.8333

$ws.Range("H122").Value = 1848.5769
$ws.Range("I122").Value = 1703.2354
$ws.Range("K122").Value = 5109.706200000001
$ws.Range("M122").Value = -2659.706200000001

$ws.Range("H132").Value = 4365.5454
$ws.Range("I132").Value = 3688.625
$ws.Range("J132").Value = 6170.6665
$ws.Range("K132").Value = 11065.875
$ws.Range("L132").Value = 18511.9995
$ws.Range("M132").Value = -8535.875
$ws.Range("N132").Value = -23571.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3423.5
$ws.Range("I82").Value = 1867.3334
$ws.Range("J82").Value = 3942.2222
$ws.Range("K82").Value = 1867.3334
$ws.Range("L82").Value = 3942.2222
$ws.Range("M82").Value = -1506.3334
$ws.Range("N82").Value = -4664.2222

$ws.Range("H85").Value = 3423.5
$ws.Range("I85").Value = 1867.3334
$ws.Range("J85").Value = 3942.2222
$ws.Range("K85").Value = 1867.3334
$ws.Range("L85").Value = 3942.2222
$ws.Range("M85").Value = -619.3334
$ws.Range("N85").Value = -6438.2222

$ws.Range("H93").Value = 998
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 1247
$ws.Range("K93").Value = 500
$ws.Range("L93").Value = 1247
$ws.Range("M93").Value = 748
$ws.Range("N93").Value = -3743

$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -2882

$ws.Range("H132").Value = 1477.2916
$ws.Range("I132").Value = 1266.8182
$ws.Range("J132").Value = 1655.3846
$ws.Range("K132").Value = 3800.4546
$ws.Range("L132").Value = 4966.1538
$ws.Range("M132").Value = -1270.4546
$ws.Range("N132").Value = -10026.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 46166.668
$ws.Range("J123").Value = 46166.668
$ws.Range("L123").Value = 46166.668
$ws.Range("N123").Value = -55966.668

$ws.Range("H126").Value = 4085.8438
$ws.Range("I126").Value = 3261.0454
$ws.Range("K126").Value = 9783.136200000001
$ws.Range("M126").Value = -7313.136200000001

$ws.Range("H132").Value = 3989.0557
$ws.Range("I132").Value = 1490.4
$ws.Range("K132").Value = 4471.200000000001
$ws.Range("M132").Value = -1941.200000000001

$ws.Range("H136").Value = 3858
$ws.Range("I136").Value = 3497.1667
$ws.Range("J136").Value = 4291
$ws.Range("K136").Value = 10491.5001
$ws.Range("L136").Value = 12873
$ws.Range("M136").Value = -7941.500100000001
$ws.Range("N136").Value = -17973
